$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "676×4=2704" "829×3=2487"
Replace-Text "724×4=2896" "402×8=3216"
Replace-Text "759×2=1518" "192×7=1344"
Replace-Text "721×7=5047" "297×8=2376"
Replace-Text "203×5=1015" "224×5=1120"
Replace-Text "994×8=7952" "124×3=372"
Replace-Text "718×8=5744" "565×9=5085"
Replace-Text "492×7=3444" "321×8=2568"
Replace-Text "434×7=3038" "663×9=5967"
Replace-Text "223×2=446" "513×9=4617"
Replace-Text "888×7=6216" "303×6=1818"
Replace-Text "498×2=996" "913×6=5478"
Replace-Text "134×6=804" "801×3=2403"
Replace-Text "999×6=5994" "261×3=783"
Replace-Text "898×9=8082" "632×8=5056"
Replace-Text "421×9=3789" "447×5=2235"
Replace-Text "946×2=1892" "146×5=730"
Replace-Text "856×3=2568" "247×3=741"
Replace-Text "219×4=876" "961×8=7688"
Replace-Text "325×6=1950" "274×5=1370"
Replace-Text "730×4=2920" "817×9=7353"
Replace-Text "714×5=3570" "223×8=1784"
Replace-Text "986×5=4930" "629×3=1887"
Replace-Text "489×5=2445" "548×4=2192"
